$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "XtGmt288"
$ws.Range("B2").Value = 23082311
$ws.Range("C2").Value = "lrdubcp55"
$ws.Range("D2").Value = "W%!Ew7p6"
$ws.Range("F2").Value = "QRsXPmJH"
$ws.Range("G2").Value = "LOuB"
